$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.10681414604187
$ws.Range("B1").Value = 1.942057013511658
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.040972709655762
$ws.Range("E1").Value = 1.115364909172058
